# ImpedanceTemplate.xlsx: extend the numbered measurement table on Sheet1
# from 64 rows (1-64) to 96 rows (1-96).
#
# Before: rows 11-74 hold sequence numbers 1-64, row 75 is a thin spacer
# row, row 76 is the bottom "Notes" bar.
# After:  rows 11-106 hold sequence numbers 1-96, the spacer row moves to
# 107 and the "Notes" bar moves to 108.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the spacer/Notes rows down by 32 rows, inserting blank rows in
# their place (shifts r75:r76 -> r107:r108).
$ws.Rows("75:106").Insert(-4121)

# The newly inserted rows have no formatting yet - stamp them with the
# same look as the rest of the numbered rows (border/fill/number style)
# by copying the last numbered row (64, row 74) across the new block.
$ws.Range("C74:J74").Copy() | Out-Null
$ws.Range("C75:J106").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Continue the 1..64 sequence in column C through 65..96.
for ($i = 0; $i -lt 32; $i++) {
    $r = 75 + $i
    $ws.Cells.Item($r, 3).Value = 65 + $i
}

# The first of the newly inserted rows picked up a slightly taller row
# height than the rest of the block.
$ws.Rows(75).RowHeight = 17.4

# Leave the selection where the edit finished.
$ws.Range("L102").Select() | Out-Null
